$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Header date
Replace-Text "2025-03-19 Wednesday" "2025-03-20 Thursday"

# Table cells - processed in document order to avoid collisions between
# an "old" value that equals a later "new" value (e.g. 74÷8=)
Replace-Text "57÷9=" "29÷8="
Replace-Text "56÷7=" "92÷8="
Replace-Text "71÷7=" "83÷5="
Replace-Text "60÷4=" "71÷2="
Replace-Text "18÷6=" "59÷6="

Replace-Text "28÷7=" "98÷2="
Replace-Text "18÷2=" "89÷6="
Replace-Text "74÷8=" "85÷4="
Replace-Text "27÷7=" "34÷2="
Replace-Text "69÷5=" "94÷8="

Replace-Text "26÷8=" "44÷3="
Replace-Text "37÷4=" "53÷5="
Replace-Text "83÷2=" "18÷8="
Replace-Text "59÷9=" "76÷9="
Replace-Text "86÷6=" "96÷9="

Replace-Text "53÷7=" "74÷8="
Replace-Text "35÷3=" "75÷8="
Replace-Text "49÷7=" "26÷3="
Replace-Text "25÷4=" "65÷3="
Replace-Text "11÷7=" "23÷3="

Replace-Text "99÷3=" "56÷3="
Replace-Text "71÷6=" "17÷4="
Replace-Text "84÷6=" "35÷4="
Replace-Text "32÷8=" "97÷5="
Replace-Text "77÷2=" "95÷2="
